$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 10:08"

# Row 7 - Rusia
$ws.Range("B7").Value = 951897
$ws.Range("C7").Value = 4921
$ws.Range("D7").Value = 767477
$ws.Range("E7").Value = 168110
$ws.Range("G7").Value = 121
$ws.Range("H7").Value = 16310

# Row 49 - Singapur
$ws.Range("B49").Value = 56266
$ws.Range("C49").Value = 50
$ws.Range("E49").Value = 2588

# Row 57 - Armenia
$ws.Range("B57").Value = 42616
$ws.Range("C57").Value = 139
$ws.Range("D57").Value = 35907
$ws.Range("E57").Value = 5859
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 850

# Row 62 - Afganistan
$ws.Range("B62").Value = 37953
$ws.Range("C62").Value = 59
$ws.Range("E62").Value = 8552

# Row 109 - Hungria
$ws.Range("B109").Value = 5133
$ws.Range("C109").Value = 35
$ws.Range("D109").Value = 3692
$ws.Range("E109").Value = 830

# Row 151 - Letonia
$ws.Range("B151").Value = 1333
$ws.Range("C151").Value = 3
$ws.Range("E151").Value = 207
